$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" message on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.76 = 51531.39 pesos`n✅ 51531.39 pesos = 12.72 = 966.53 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate cells on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 78.36
$wsTasas.Range("O10").Value = 4038
$wsTasas.Range("N12").Value = 4052
$wsTasas.Range("O12").Value = 76
